$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 (Chuk Yuen (South) Estate): set Visit Date and Progress
$ws.Range("L20").Value = 43978
$ws.Range("L20").NumberFormat = "d-mmm"
$ws.Range("M20").Value = "Webpage Done"

# Row 21 (Lok Fu Estate): add hyperlink to the URL already shown in K21,
# and set Visit Date and Progress
$ws.Hyperlinks.Add($ws.Range("K21"), "https://www.housingauthority.gov.hk/en/global-elements/estate-locator/detail.html?propertyType=1&id=2727")
$ws.Range("L21").Value = 43978
$ws.Range("L21").NumberFormat = "d-mmm"
$ws.Range("M21").Value = "Webpage Done"

# Row 22 (Lower Wong Tai Sin (II) Estate): set Visit Date and Progress
$ws.Range("L22").Value = 43978
$ws.Range("L22").NumberFormat = "d-mmm"
$ws.Range("M22").Value = "Webpage Done"

# Update the saved selection to M24
$ws.Range("M24").Select()
